$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.838.65"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "'3.255.47"
$ws.Range("E3").Value = "  -3.20%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'578.63"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "'173.21"
$ws.Range("E6").Value = "  -7.99%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -4.08%  "
$ws.Range("D9").Value = "'3.250.66"
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("D10").Value = "'0.172"
$ws.Range("E10").Value = "  -6.93%  "
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("D12").Value = "'44.70"
$ws.Range("E12").Value = "  -6.14%  "
$ws.Range("D13").Value = "'0.0000269"
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").Value = "'665.99"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "'3.775.27"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("E16").Value = "  -4.74%  "
$ws.Range("D17").Value = "'66.717.46"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D19").Value = "'3.244.52"
$ws.Range("E19").Value = "  -3.61%  "
$ws.Range("D20").Value = "'17.20"
$ws.Range("E20").Value = "  -4.76%  "
$ws.Range("D21").Value = "'10.72"
$ws.Range("E21").Value = "  -4.40%  "
$ws.Range("E22").Value = "  -3.79%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "'5.29"
$ws.Range("E23").Value = "  +3.63%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'16.88"
$ws.Range("E24").Value = "  -6.78%  "
$ws.Range("D25").Value = "'96.54"
$ws.Range("E25").Value = "  -4.45%  "
$ws.Range("D26").Value = "'3.84"
$ws.Range("E26").Value = "  -4.62%  "
$ws.Range("D27").Value = "'2.63"
$ws.Range("E27").Value = "  -7.68%  "
$ws.Range("D28").Value = "'9.10"
$ws.Range("E28").Value = "  -7.13%  "
$ws.Range("D29").Value = "'32.02"
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("E30").Value = "  -4.92%  "
$ws.Range("D31").Value = "'6.81"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "'564.42"
$ws.Range("E32").Value = "  -8.43%  "
$ws.Range("D33").Value = "'10.86"
$ws.Range("E33").Value = "  -3.24%  "
$ws.Range("D34").Value = "'3.749.37"
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("E35").Value = "  -4.41%  "
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").Value = "'3.44"
$ws.Range("E36").Value = "  -12.60%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'55.44"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("E40").Value = "  -4.72%  "
$ws.Range("D41").Value = "'2.60"
$ws.Range("E41").Value = "  -8.15%  "
$ws.Range("D42").Value = "'0.0₃0656"
$ws.Range("E42").Value = "  -7.23%  "
$ws.Range("D43").Value = "'3.00"
$ws.Range("E43").Value = "  -8.11%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.324"
$ws.Range("E44").Value = "  -6.42%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.21"
$ws.Range("E45").Value = "  -5.44%  "
$ws.Range("E46").Value = "  -6.54%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("E51").Value = "  -0.88%  "
